$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Add($ws.Range("D23"), "mailto:samsungbraasilpark.telecel@gmail.com")
Write-Host "done"
